# "Generate Report for Archive"
# The status "Ready for handoff" has moved on to "In Translation" for the
# files tracked in this handback report, and the Status column is
# re-sized (narrower, since the new text is shorter than the old one).

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"
# Target stored column width is 13.4101845877511 characters. The host's
# column-width grid only lands on multiples of 1/6, so feed it the
# ColumnWidth input whose post-quantization value (13.333333333333334)
# is the closest representable match to that target.
$newColWidth = 12.501302083333332

# --- Overview sheet: Status is mirrored into columns E (zh-cn) and F (de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in 2..4) {
    foreach ($col in @("E", "F")) {
        $cell = $wsOverview.Range("$col$r")
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}
$wsOverview.Columns.Item(5).ColumnWidth = $newColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColWidth

# --- Per-locale detail sheets: Status lives in column C
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($r in 2..4) {
        $cell = $ws.Range("C$r")
        if ($cell.Value2 -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
    $ws.Columns.Item(3).ColumnWidth = $newColWidth
}
